$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 11: 2016-11-18, DEV, "Start omzetten html5 template naar drupal template", 2u ---
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A11").Value = 42692
$ws.Range("B11").Value = "DEV "
$ws.Range("C11").Value = "Start omzetten html5 template naar drupal template"
$ws.Range("D11").Value = "2u"

# --- New row 12: 2016-11-20, DES + DEV, "Afwerking drupal template + dossier aanpassen", 3u ---
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = 42694
$ws.Range("B12").Value = "DES + DEV"
$ws.Range("C12").Value = "Afwerking drupal template + dossier aanpassen"
$ws.Range("D12").Value = "3u"

# --- Update the running total in F2 ---
$ws.Range("F2").Value = "Totaal: 28u"

# --- Update the selection to match the new active cell ---
$ws.Range("F18").Select() | Out-Null

Write-Host "done"
